$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 156863
$ws.Range("C4").Value = 147935
$ws.Range("C5").Value = 8928
$ws.Range("C8").Value = 63.96
